# Updates cryptos list values to match the latest scrape (GitHub Actions run).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force each touched cell to remain plain Text so numeric-looking strings
# (e.g. "0.450", "0.0000237") are not silently re-typed as numbers and
# lose their exact on-disk representation, matching the original inlineStr cells.
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "63.292.22"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "  -0.02%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.056.44"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "  -0.48%  "
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "  +0.01%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "591.55"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "  +0.50%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "153.03"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "  +0.72%  "
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "  -0.05%  "
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "  -2.06%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "3.056.76"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "  -0.26%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.154"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "  -0.79%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "5.81"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "  -1.22%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.450"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "  -2.55%  "
$ws.Range("B13").NumberFormat = "@"
$ws.Range("B13").Value = "ShibaInu"
$ws.Range("C13").NumberFormat = "@"
$ws.Range("C13").Value = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.0000237"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "  -1.94%  "
$ws.Range("B14").NumberFormat = "@"
$ws.Range("B14").Value = "Avalanche"
$ws.Range("C14").NumberFormat = "@"
$ws.Range("C14").Value = "https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "36.58"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "  -1.61%  "
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "  +1.59%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "3.561.34"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "  -0.53%  "
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "  -0.53%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "63.278.13"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "  -0.12%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "3.059.16"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "  -0.36%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "482.49"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "  +1.05%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "14.31"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "  -2.03%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.707"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "  -1.49%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "7.53"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "  +0.10%  "
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "  +3.14%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "82.14"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "  +0.57%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "12.87"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "  -1.47%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "10.75"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "  +10.25%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0.999"
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "  +0.00%  "
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = "  +1.21%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "2.68"
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = "  +0.37%  "
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = "  +0.97%  "
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = "  +0.11%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "27.47"
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = "  +0.56%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.111"
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = "  -2.14%  "
$ws.Range("B35").NumberFormat = "@"
$ws.Range("B35").Value = "Mantle"
$ws.Range("C35").NumberFormat = "@"
$ws.Range("C35").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.06"
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = "  +1.13%  "
$ws.Range("B36").NumberFormat = "@"
$ws.Range("B36").Value = "PEPE"
$ws.Range("C36").NumberFormat = "@"
$ws.Range("C36").Value = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.0₃0820"
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = "  -3.85%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "3.32"
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = "  -1.32%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "5.96"
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "  -2.64%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "2.23"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "  +0.99%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "9.28"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "  -0.96%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "50.58"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "438.16"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "  -1.29%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.290"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "  +1.49%  "
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "  +3.79%  "
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "  +0.14%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "2.833.72"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "  +0.94%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "38.36"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "  -4.20%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "129.63"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "  -0.88%  "
$ws.Range("B49").NumberFormat = "@"
$ws.Range("B49").Value = "USDe"
$ws.Range("C49").NumberFormat = "@"
$ws.Range("C49").Value = "https://coinranking.com/coin/exbfr2U-0+usde-usde"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.999"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "  +0.00%  "
$ws.Range("B50").NumberFormat = "@"
$ws.Range("B50").Value = "InjectiveProtocol"
$ws.Range("C50").NumberFormat = "@"
$ws.Range("C50").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "25.27"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "  +0.76%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "2.23"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "  -1.31%  "
